# Reclassify a handful of Choice answers in electric_en_components so the
# questions land in more homogeneous groups (see commit message).
#
#   d) Ground -> d) Earth            (Altern current source / Battery / Motor / Ground)
#   a) Botón  -> a) Rele             (Botón / Pushbutton / Switch / Two way Switch)
#   c) Ground -> c) Earth            (Battery / Receiver / Ground / Antenna)
#   a) Mass   -> a) Ground plane     (Mass / Battery / Ground / Brush)
#   c) Ground -> c) Earth            (   "          "      "      "   )
#   a) Ground -> a) Ground plane     (Ground / Mass / Connector / ...)

$d = $word.ActiveDocument

function Get-ParaText($para) {
    if ($para -eq $null) { return "" }
    # Paragraph.Range.Text carries a trailing paragraph-mark (CR) character;
    # strip it so comparisons are against the visible text only.
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

function Replace-ChoiceText($para, $oldText, $newText) {
    # Narrow the edit to just the option text itself (after the "x)" + tab
    # marker) so the rest of the run is left completely untouched.
    $full = $para.Range
    $fullText = $full.Text
    $relIdx = $fullText.IndexOf($oldText)
    if ($relIdx -lt 0) {
        throw "Expected text '$oldText' not found in paragraph text '$fullText'"
    }
    $absStart = $full.Start + $relIdx
    $absEnd = $absStart + $oldText.Length
    $target = $d.Range($absStart, $absEnd)
    $target.Text = $newText
}

$count = $d.Paragraphs.Count
$editsMade = 0

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $cur = Get-ParaText $p
    $prev1 = Get-ParaText ($(if ($i -gt 1) { $d.Paragraphs.Item($i - 1) } else { $null }))
    $next1 = Get-ParaText ($(if ($i -lt $count) { $d.Paragraphs.Item($i + 1) } else { $null }))

    if ($cur -eq "d)`tGround" -and $prev1 -eq "c)`tMotor") {
        Replace-ChoiceText $p "Ground" "Earth"
        $editsMade++
    }
    elseif ($cur -eq "a)`tBotón") {
        Replace-ChoiceText $p "Botón" "Rele"
        $editsMade++
    }
    elseif ($cur -eq "c)`tGround" -and $prev1 -eq "b)`tReceiver") {
        Replace-ChoiceText $p "Ground" "Earth"
        $editsMade++
    }
    elseif ($cur -eq "a)`tMass" -and $next1 -eq "b)`tBattery") {
        Replace-ChoiceText $p "Mass" "Ground plane"
        $editsMade++
    }
    elseif ($cur -eq "c)`tGround" -and $prev1 -eq "b)`tBattery") {
        Replace-ChoiceText $p "Ground" "Earth"
        $editsMade++
    }
    elseif ($cur -eq "a)`tGround" -and $next1 -eq "b)`tMass") {
        Replace-ChoiceText $p "Ground" "Ground plane"
        $editsMade++
    }
}

if ($editsMade -ne 6) {
    throw "Expected to make 6 choice-text edits, but made $editsMade"
}

Write-Host "Done applying choice text edits ($editsMade)."
